# Improve bathtub matching / add a second Walls test case (WLL-003)
# Target: Walls worksheet gains row 4 (A4:L4) and its used range grows
# from A1:L3 to A1:L4.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Walls")

# New row of test data describing a second Alcove Tub Wall Kit (Maax)
$ws.Range("A4").Value = "WLL-003"
$ws.Range("B4").Value = "Test Alcove Tub Wall Kit Maax"

# C4 stays blank (mirrors the blank Image URL cells in C2/C3); copy an
# existing blank cell so the cell is materialized instead of being left
# out of the sheet entirely.
$ws.Range("C2").Copy($ws.Range("C4"))

$ws.Range("D4").Value = "60 x 32"
$ws.Range("E4").Value = 60
$ws.Range("F4").Value = 32
$ws.Range("G4").Value = "Alcove Tub Wall Kit"
$ws.Range("H4").Value = "Yes"
$ws.Range("I4").Value = "Maax"
$ws.Range("J4").Value = "MAAX"
$ws.Range("K4").Value = "Utile"
$ws.Range("L4").Value = 950
